$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the modified price/volume cells keep a Text number format so that
# values such as "327.03" and "-1.08%" are stored as literal text, matching
# the original inline-string cells rather than being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "327.03"
$ws.Range("E2").Value = "-1.08%"
$ws.Range("D3").Value = "43.53"
$ws.Range("E3").Value = "5.01%"
$ws.Range("D4").Value = "5.538"
$ws.Range("E4").Value = "-2.79%"
$ws.Range("D5").Value = "0.08084"
$ws.Range("E5").Value = "-4.08%"
$ws.Range("D6").Value = "8.641"
$ws.Range("E6").Value = "-2.00%"
$ws.Range("D7").Value = "4.286"
$ws.Range("E7").Value = "-4.69%"
$ws.Range("D8").Value = "1.890"
$ws.Range("E8").Value = "-4.78%"
$ws.Range("E9").Value = "-6.18%"
$ws.Range("D10").Value = "0.9383"
$ws.Range("E10").Value = "1.22%"
$ws.Range("D11").Value = "0.1182"
$ws.Range("E11").Value = "-6.32%"
$ws.Range("D12").Value = "0.1888"
$ws.Range("E12").Value = "-4.80%"
$ws.Range("D13").Value = "0.09575"
$ws.Range("E13").Value = "0.67%"
$ws.Range("D14").Value = "0.04129"
$ws.Range("E14").Value = "4.22%"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").Value = "0.33%"
$ws.Range("D16").Value = "0.001274"
$ws.Range("E16").Value = "-2.21%"
$ws.Range("D17").Value = "0.005911"
$ws.Range("E17").Value = "-3.31%"
$ws.Range("E18").Value = "3.73%"
$ws.Range("D19").Value = "0.3486"
$ws.Range("E19").Value = "-0.72%"
$ws.Range("D20").Value = "8.572"
$ws.Range("E20").Value = "-6.50%"
$ws.Range("D21").Value = "0.1362"
$ws.Range("E21").Value = "-0.16%"
$ws.Range("D22").Value = "0.2589"
$ws.Range("E22").Value = "3.09%"
$ws.Range("D23").Value = "0.04327"
$ws.Range("E23").Value = "-1.72%"
$ws.Range("D24").Value = "0.001240"
$ws.Range("E24").Value = "-0.51%"
$ws.Range("D25").Value = "0.004384"
$ws.Range("E25").Value = "-0.32%"
$ws.Range("E26").Value = "3.39%"
$ws.Range("D27").Value = "0.0003995"
$ws.Range("E27").Value = "-0.01%"
$ws.Range("D39").Value = "0.02646"
$ws.Range("E39").Value = "-6.52%"
$ws.Range("D40").Value = "0.05457"
$ws.Range("E40").Value = "-1.16%"
$ws.Range("D41").Value = "0.01145"
$ws.Range("E41").Value = "27.86%"
$ws.Range("D42").Value = "0.007720"
$ws.Range("E42").Value = "-2.37%"
$ws.Range("D43").Value = "0.1395"
$ws.Range("E43").Value = "-3.01%"
$ws.Range("D45").Value = "0.009245"
$ws.Range("D46").Value = "0.00007008"
$ws.Range("E46").Value = "-4.36%"
$ws.Range("E47").Value = "0.00%"
$ws.Range("D48").Value = "0.003564"
$ws.Range("E48").Value = "10.64%"
$ws.Range("D49").Value = "0.002273"
$ws.Range("E49").Value = "-0.34%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.00%"
